$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9487097263336182
$ws.Range("B1").Value = 3.672916173934937
$ws.Range("C1").Value = 2.938079595565796
$ws.Range("D1").Value = 2.541958093643188
$ws.Range("E1").Value = 2.068214654922485
